$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: push the existing rows 21 and 22 (current "Sweet Heart" / "Brooks" 10kg
# entries) down to new rows 23 and 24, preserving values and number formats.
# Column D (4) carries a date number format that must be copied explicitly;
# the other columns use the default "General" style so it is left untouched.
$ws.Cells.Item(23, 4).NumberFormat = $ws.Cells.Item(21, 4).NumberFormat
$ws.Cells.Item(24, 4).NumberFormat = $ws.Cells.Item(22, 4).NumberFormat

for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(23, $col).Value = $ws.Cells.Item(21, $col).Value()
    $ws.Cells.Item(24, $col).Value = $ws.Cells.Item(22, $col).Value()
}

# Step 2: overwrite row 21 with the new weekly "Lapins" (10 kilos bandeja) entry.
$ws.Range("D21").Value = 44931
$ws.Range("K21").Value = "Lapins"
$ws.Range("M21").Value = 250
$ws.Range("N21").Value = 6000
$ws.Range("O21").Value = 6500
$ws.Range("P21").Value = 6250
$ws.Range("S21").Value = 625

# Step 3: overwrite row 22 with the new weekly "Lapins" (5 kilos bandeja) entry.
$ws.Range("D22").Value = 44931
$ws.Range("K22").Value = "Lapins"
$ws.Range("N22").Value = 3000
$ws.Range("O22").Value = 3300
$ws.Range("P22").Value = 3150
$ws.Range("Q22").Value = "$/bandeja 5 kilos"
$ws.Range("S22").Value = 630
$ws.Range("T22").Value = 5
